# Replace "Agri Fund" with "SAAS Fund" throughout the fund column,
# and move the active selection to A4 (as captured by the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SAAS Fund"
$ws.Range("A3").Value = "SAAS Fund"
$ws.Range("A4").Value = "SAAS Fund"

$ws.Range("A4").Select()
